# Acta de Conciliación template: correct the merge-field placeholder
# used for the conciliation proposal(s) paragraph.
#
#   {propuestas}  ->  {propuesta}
#
# (The rest of the upstream diff is just Word's proofing engine
#  re-splitting existing runs and inserting <w:proofErr/> markers around
#  them after a spell/grammar check pass - the visible text is unchanged
#  everywhere else, so no further edits are required.)

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "{propuestas}",  # FindText
    $true,           # MatchCase
    $true,           # MatchWholeWord
    $false,          # MatchWildcards
    $false,          # MatchSoundsLike
    $false,          # MatchAllWordForms
    $true,           # Forward
    1,               # Wrap            (wdFindContinue)
    $false,          # Format
    "{propuesta}",   # ReplaceWith
    2                # Replace         (wdReplaceAll)
)
